{"js": "const replacements = [\n  [\"2025-06-20 Friday\", \"2025-06-21 Saturday\"],\n  [\"121\u00f74=30, 1\", \"138\u00f79=15, 3\"],\n  [\"384\u00f75=76, 4\", \"470\u00f74=117, 2\"],\n  [\"765\u00f76=127, 3\", \"948\u00f76=158, 0\"],\n  [\"139\u00f77=19, 6\", \"595\u00f77=85, 0\"],\n  [\"809\u00f76=134, 5\", \"158\u00f77=22, 4\"],\n  [\"756\u00f78=94, 4\", \"989\u00f72=494, 1\"],\n  [\"607\u00f78=75, 7\", \"545\u00f78=68, 1\"],\n  [\"307\u00f72=153, 1\", \"146\u00f75=29, 1\"],\n  [\"990\u00f73=330, 0\", \"977\u00f72=488, 1\"],\n  [\"882\u00f79=98, 0\", \"200\u00f73=66, 2\"],\n  [\"446\u00f75=89, 1\", \"219\u00f73=73, 0\"],\n  [\"428\u00f72=214, 0\", \"661\u00f76=110, 1\"],\n  [\"928\u00f79=103, 1\", \"291\u00f73=97, 0\"],\n  [\"465\u00f73=155, 0\", \"689\u00f73=229, 2\"],\n  [\"455\u00f79=50, 5\", \"770\u00f72=385, 0\"],\n  [\"199\u00f72=99, 1\", \"382\u00f79=42, 4\"],\n  [\"485\u00f77=69, 2\", \"912\u00f75=182, 2\"],\n  [\"706\u00f76=117, 4\", \"302\u00f77=43, 1\"],\n  [\"521\u00f74=130, 1\", \"726\u00f75=145, 1\"],\n  [\"351\u00f75=70, 1\", \"425\u00f78=53, 1\"],\n  [\"926\u00f79=102, 8\", \"659\u00f79=73, 2\"],\n  [\"487\u00f73=162, 1\", \"910\u00f72=455, 0\"],\n  [\"264\u00f76=44, 0\", \"998\u00f74=249, 2\"],\n  [\"429\u00f74=107, 1\", \"144\u00f75=28, 4\"],\n  [\"193\u00f77=27, 4\", \"932\u00f79=103, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2025-06-20 Friday\", \"2025-06-21 Saturday\")\n    ,@(\"121\u00f74=30, 1\", \"138\u00f79=15, 3\")\n    ,@(\"384\u00f75=76, 4\", \"470\u00f74=117, 2\")\n    ,@(\"765\u00f76=127, 3\", \"948\u00f76=158, 0\")\n    ,@(\"139\u00f77=19, 6\", \"595\u00f77=85, 0\")\n    ,@(\"809\u00f76=134, 5\", \"158\u00f77=22, 4\")\n    ,@(\"756\u00f78=94, 4\", \"989\u00f72=494, 1\")\n    ,@(\"607\u00f78=75, 7\", \"545\u00f78=68, 1\")\n    ,@(\"307\u00f72=153, 1\", \"146\u00f75=29, 1\")\n    ,@(\"990\u00f73=330, 0\", \"977\u00f72=488, 1\")\n    ,@(\"882\u00f79=98, 0\", \"200\u00f73=66, 2\")\n    ,@(\"446\u00f75=89, 1\", \"219\u00f73=73, 0\")\n    ,@(\"428\u00f72=214, 0\", \"661\u00f76=110, 1\")\n    ,@(\"928\u00f79=103, 1\", \"291\u00f73=97, 0\")\n    ,@(\"465\u00f73=155, 0\", \"689\u00f73=229, 2\")\n    ,@(\"455\u00f79=50, 5\", \"770\u00f72=385, 0\")\n    ,@(\"199\u00f72=99, 1\", \"382\u00f79=42, 4\")\n    ,@(\"485\u00f77=69, 2\", \"912\u00f75=182, 2\")\n    ,@(\"706\u00f76=117, 4\", \"302\u00f77=43, 1\")\n    ,@(\"521\u00f74=130, 1\", \"726\u00f75=145, 1\")\n    ,@(\"351\u00f75=70, 1\", \"425\u00f78=53, 1\")\n    ,@(\"926\u00f79=102, 8\", \"659\u00f79=73, 2\")\n    ,@(\"487\u00f73=162, 1\", \"910\u00f72=455, 0\")\n    ,@(\"264\u00f76=44, 0\", \"998\u00f74=249, 2\")\n    ,@(\"429\u00f74=107, 1\", \"144\u00f75=28, 4\")\n    ,@(\"193\u00f77=27, 4\", \"932\u00f79=103, 5\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
